# Tactile Tabletop Data - Augments workbook update
# "adjusting for more buffs, and creating accurate and current printing document"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tactile Tabletop Data - Level 1")

# --- Update the "Conditions to Remove" text on the existing transforming
#     augments to note that the card transforms afterwards ---
$ws.Range("D2").Value = "Roll a 3 four times, then card transforms"
$ws.Range("D3").Value = "Wait a week, then card transforms"
$ws.Range("D4").Value = "Drink a gallon of blood, then card transforms"
$ws.Range("D5").Value = "Heal 50 points worth of health from rest, then card transforms"
$ws.Range("D6").Value = "Become downed three times, then card transforms"

# --- Row 8 ("Spirit Boon") is replaced by a new augment, "Spectre of Death" ---
$ws.Range("A8").Value = "Spectre of Death"
$ws.Range("B8").Value = "Attack - discard"
$ws.Range("C8").Value = "A mostly transparent spectre floats just behind you at all times, about an inch behind and above you. At the start of each turn, choose to increase your influence, defense, or attack by 1"
$ws.Rows.Item(8).RowHeight = 90

# --- New augment rows 9-11 ---
$ws.Range("A9").Value = "Eyes of the fallen"
$ws.Range("B9").Value = "Attack - discard"
$ws.Range("C9").Value = "Your eyes are blackened completely. You gain +4 to any roll that involves percieving "
$ws.Rows.Item(9).RowHeight = 45

$ws.Range("A10").Value = "Flames of ambition"
$ws.Range("B10").Value = "Attack - discard"
$ws.Range("C10").Value = "Half of your body is scarred. Gain +2 to influence values when within 15 ft of an open flame"
$ws.Rows.Item(10).RowHeight = 45

$ws.Range("A11").Value = "Beastial Mark"
$ws.Range("B11").Value = "Attack - discard"
$ws.Range("C11").Value = "Small horns grow from the top of your head, and your body hair grows thicker. Defense values are increased by 2"
$ws.Rows.Item(11).RowHeight = 60

# --- Remove the stray formatted-but-empty cell at C19 (row no longer used) ---
$ws.Rows.Item(19).Delete()

# --- Update the remembered selection on this sheet ---
$ws.Range("F5").Select() | Out-Null
